$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.513.72'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.917.95'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = '  +0.80%  '
$ws.Range("D5").Value = "'325.37"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").Value = "'0.4811"
$ws.Range("E7").Value = '  -0.75%  '
$ws.Range("D8").Value = "'0.4050"
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").Value = "'0.08214"
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("D11").Value = "'23.42"
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").Value = '1.909.06'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = "'6.049"
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = "'7.231"
$ws.Range("E14").Value = '  +1.68%  '
$ws.Range("D15").Value = "'91.44"
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("D16").Value = "'0.06881"
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = "'17.52"
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").Value = '29.518.82'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '2.131.71'
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("D26").Value = "'6.534"
$ws.Range("E26").Value = '  +3.33%  '
$ws.Range("D27").Value = "'155.91"
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("D28").Value = "'20.02"
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("D29").Value = "'2.096"
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").Value = "'120.59"
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").Value = "'1.016"
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("D32").Value = "'0.09635"
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").Value = "'5.615"
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("D34").Value = "'3.559"
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").Value = "'1.372"
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").Value = "'0.06305"
$ws.Range("E36").Value = '  +3.00%  '
$ws.Range("D37").Value = "'0.02283"
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("D38").Value = "'1.186"
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").Value = "'0.5933"
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").Value = "'10.69"
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").Value = "'7.909"
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").Value = "'0.1846"
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").Value = "'2.464"
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("D44").Value = "'1.280"
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = "'12.34"
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").Value = "'0.07468"
$ws.Range("E46").Value = '  -3.16%  '
$ws.Range("D47").Value = "'0.5557"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").Value = "'1.939"
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("D49").Value = "'118.19"
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("D50").Value = "'2.427"
$ws.Range("E50").Value = '  +3.20%  '
$ws.Range("D51").Value = "'71.96"
$ws.Range("E51").Value = '  -1.17%  '
